$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam1"
$ws.Range("C2").Value = "Spn"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.10268633333333
$ws.Range("H2").Value = 60.308059
$ws.Range("I2").Value = 0.1188668172183431
$ws.Range("J2").Value = 0.1192185838730403
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.359683333333334
$ws.Range("N2").Value = 16.07905
$ws.Range("O2").Value = 0.5495559766256753
$ws.Range("P2").Value = 0.5567946816040513
$ws.Range("Q2").Value = 107.7440328959945
$ws.Range("R2").Value = 969.6962960639502
$ws.Range("S2").Value = 0.06532396982481216
$ws.Range("T2").Value = 0.06638027344887537

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam1"
$ws.Range("C3").Value = "Spn"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.10268633333333
$ws.Range("H3").Value = 60.308059
$ws.Range("I3").Value = 0.1188668172183431
$ws.Range("J3").Value = 0.1192185838730403
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.380377
$ws.Range("N3").Value = 0.760754
$ws.Range("O3").Value = 0.03900201573866823
$ws.Range("P3").Value = 0.02634383133387908
$ws.Range("Q3").Value = 7.646599519414334
$ws.Range("R3").Value = 45.879597116486
$ws.Range("S3").Value = 0.004636045475955216
$ws.Range("T3").Value = 0.00314067426541529

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Icam1"
$ws.Range("C4").Value = "Spn"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.10268633333333
$ws.Range("H4").Value = 60.308059
$ws.Range("I4").Value = 0.1188668172183431
$ws.Range("J4").Value = 0.1192185838730403
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.012692
$ws.Range("N4").Value = 12.038076
$ws.Range("O4").Value = 0.4114420076356565
$ws.Range("P4").Value = 0.4168614870620697
$ws.Range("Q4").Value = 80.66588862827601
$ws.Range("R4").Value = 725.992997654484
$ws.Range("S4").Value = 0.0489068019175757
$ws.Range("T4").Value = 0.04969763615874966

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam1"
$ws.Range("C5").Value = "Spn"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 29.50180766666667
$ws.Range("H5").Value = 88.50542300000001
$ws.Range("I5").Value = 0.1744436500364427
$ws.Range("J5").Value = 0.1749598871212952
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.359683333333334
$ws.Range("N5").Value = 16.07905
$ws.Range("O5").Value = 0.5495559766256753
$ws.Range("P5").Value = 0.5567946816040513
$ws.Range("Q5").Value = 158.1203468542389
$ws.Range("R5").Value = 1423.08312168815
$ws.Range("S5").Value = 0.09586655046192477
$ws.Range("T5").Value = 0.09741673464318232

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Icam1"
$ws.Range("C6").Value = "Spn"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 29.50180766666667
$ws.Range("H6").Value = 88.50542300000001
$ws.Range("I6").Value = 0.1744436500364427
$ws.Range("J6").Value = 0.1749598871212952
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.380377
$ws.Range("N6").Value = 0.760754
$ws.Range("O6").Value = 0.03900201573866823
$ws.Range("P6").Value = 0.02634383133387908
$ws.Range("Q6").Value = 11.22180909482367
$ws.Range("R6").Value = 67.33085456894202
$ws.Range("S6").Value = 0.00680365398423207
$ws.Range("T6").Value = 0.004609113756517924

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Icam1"
$ws.Range("C7").Value = "Spn"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 29.50180766666667
$ws.Range("H7").Value = 88.50542300000001
$ws.Range("I7").Value = 0.1744436500364427
$ws.Range("J7").Value = 0.1749598871212952
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.012692
$ws.Range("N7").Value = 12.038076
$ws.Range("O7").Value = 0.4114420076356565
$ws.Range("P7").Value = 0.4168614870620697
$ws.Range("Q7").Value = 118.381667609572
$ws.Range("R7").Value = 1065.435008486148
$ws.Range("S7").Value = 0.07177344559028584
$ws.Range("T7").Value = 0.07293403872159497

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Icam1"
$ws.Range("C8").Value = "Spn"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 61.77435033333333
$ws.Range("H8").Value = 185.323051
$ws.Range("I8").Value = 0.3652706055348701
$ws.Range("J8").Value = 0.3663515633831165
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.359683333333334
$ws.Range("N8").Value = 16.07905
$ws.Range("O8").Value = 0.5495559766256753
$ws.Range("P8").Value = 0.5567946816040513
$ws.Range("Q8").Value = 331.0909559090611
$ws.Range("R8").Value = 2979.81860318155
$ws.Range("S8").Value = 0.2007366443573673
$ws.Range("T8").Value = 0.2039826020890488

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Icam1"
$ws.Range("C9").Value = "Spn"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 61.77435033333333
$ws.Range("H9").Value = 185.323051
$ws.Range("I9").Value = 0.3652706055348701
$ws.Range("J9").Value = 0.3663515633831165
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.380377
$ws.Range("N9").Value = 0.760754
$ws.Range("O9").Value = 0.03900201573866823
$ws.Range("P9").Value = 0.02634383133387908
$ws.Range("Q9").Value = 23.49754205674233
$ws.Range("R9").Value = 140.985252340454
$ws.Range("S9").Value = 0.01424628990594388
$ws.Range("T9").Value = 0.00965110379466773

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Icam1"
$ws.Range("C10").Value = "Spn"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 61.77435033333333
$ws.Range("H10").Value = 185.323051
$ws.Range("I10").Value = 0.3652706055348701
$ws.Range("J10").Value = 0.3663515633831165
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.012692
$ws.Range("N10").Value = 12.038076
$ws.Range("O10").Value = 0.4114420076356565
$ws.Range("P10").Value = 0.4168614870620697
$ws.Range("Q10").Value = 247.881441387764
$ws.Range("R10").Value = 2230.932972489876
$ws.Range("S10").Value = 0.1502876712715589
$ws.Range("T10").Value = 0.1527178574994

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Icam1"
$ws.Range("C11").Value = "Spn"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.4970125
$ws.Range("H11").Value = 2.994025
$ws.Range("I11").Value = 0.008851807577379077
$ws.Range("J11").Value = 0.005918668690373198
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.359683333333334
$ws.Range("N11").Value = 16.07905
$ws.Range("O11").Value = 0.5495559766256753
$ws.Range("P11").Value = 0.5567946816040513
$ws.Range("Q11").Value = 8.023512946041668
$ws.Range("R11").Value = 48.14107767625001
$ws.Range("S11").Value = 0.004864563758089112
$ws.Range("T11").Value = 0.003295483248976212

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Icam1"
$ws.Range("C12").Value = "Spn"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.4970125
$ws.Range("H12").Value = 2.994025
$ws.Range("I12").Value = 0.008851807577379077
$ws.Range("J12").Value = 0.005918668690373198
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.380377
$ws.Range("N12").Value = 0.760754
$ws.Range("O12").Value = 0.03900201573866823
$ws.Range("P12").Value = 0.02634383133387908
$ws.Range("Q12").Value = 0.5694291237125001
$ws.Range("R12").Value = 2.27771649485
$ws.Range("S12").Value = 0.0003452383384486014
$ws.Range("T12").Value = 0.0001559204097003025

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Icam1"
$ws.Range("C13").Value = "Spn"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.4970125
$ws.Range("H13").Value = 2.994025
$ws.Range("I13").Value = 0.008851807577379077
$ws.Range("J13").Value = 0.005918668690373198
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.012692
$ws.Range("N13").Value = 12.038076
$ws.Range("O13").Value = 0.4114420076356565
$ws.Range("P13").Value = 0.4168614870620697
$ws.Range("Q13").Value = 6.007050082650001
$ws.Range("R13").Value = 36.0423004959
$ws.Range("S13").Value = 0.003642005480841364
$ws.Range("T13").Value = 0.002467265031696684

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Icam1"
$ws.Range("C14").Value = "Spn"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 56.24355599999999
$ws.Range("H14").Value = 168.730668
$ws.Range("I14").Value = 0.3325671196329652
$ws.Range("J14").Value = 0.3335512969321748
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.359683333333334
$ws.Range("N14").Value = 16.07905
$ws.Range("O14").Value = 0.5495559766256753
$ws.Range("P14").Value = 0.5567946816040513
$ws.Range("Q14").Value = 301.4476497006
$ws.Range("R14").Value = 2713.0288473054
$ws.Range("S14").Value = 0.182764248223482
$ws.Range("T14").Value = 0.1857195881739686

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Icam1"
$ws.Range("C15").Value = "Spn"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 56.24355599999999
$ws.Range("H15").Value = 168.730668
$ws.Range("I15").Value = 0.3325671196329652
$ws.Range("J15").Value = 0.3335512969321748
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.5
$ws.Range("M15").Value = 0.380377
$ws.Range("N15").Value = 0.760754
$ws.Range("O15").Value = 0.03900201573866823
$ws.Range("P15").Value = 0.02634383133387908
$ws.Range("Q15").Value = 21.393755100612
$ws.Range("R15").Value = 128.362530603672
$ws.Range("S15").Value = 0.01297078803408847
$ws.Range("T15").Value = 0.00878701910757783

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Icam1"
$ws.Range("C16").Value = "Spn"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 56.24355599999999
$ws.Range("H16").Value = 168.730668
$ws.Range("I16").Value = 0.3325671196329652
$ws.Range("J16").Value = 0.3335512969321748
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 4.012692
$ws.Range("N16").Value = 12.038076
$ws.Range("O16").Value = 0.4114420076356565
$ws.Range("P16").Value = 0.4168614870620697
$ws.Range("Q16").Value = 225.688067212752
$ws.Range("R16").Value = 2031.192604914768
$ws.Range("S16").Value = 0.1368320833753948
$ws.Range("T16").Value = 0.1390446896506283

Write-Output "done"